# Insert a new weekly record at row 76 (Terminal Hortofrutícola Agro Chillán - Mango),
# shifting the existing rows 76:110 down to 77:111.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Cells.Item(76, 1).Value = 7
$ws.Cells.Item(76, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(76, 3).Value = "Ñuble"
$ws.Cells.Item(76, 4).Value = 45009
$ws.Cells.Item(76, 5).Value = 16
$ws.Cells.Item(76, 6).Value = "Fruta"
$ws.Cells.Item(76, 7).Value = 100108
$ws.Cells.Item(76, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(76, 9).Value = 100108002
$ws.Cells.Item(76, 10).Value = "Mango"
$ws.Cells.Item(76, 11).Value = "Sin especificar"
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 50
$ws.Cells.Item(76, 14).Value = 7000
$ws.Cells.Item(76, 15).Value = 7000
$ws.Cells.Item(76, 16).Value = 7000
$ws.Cells.Item(76, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(76, 18).Value = "Perú"
$ws.Cells.Item(76, 19).Value = 1750
$ws.Cells.Item(76, 20).Value = 4
